$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: wrap a WordprocessingML body fragment into a full OPC "flat xml"
# package so that Range.InsertXML can consume it.
# ---------------------------------------------------------------------------
function New-BodyXml([string]$bodyFragment) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">$bodyFragment</w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gains two trailing
#    spaces and three new red (C00000) runs forming
#    "(This is a change – Version for branch alternate)"
#    split across three <w:r> elements.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)

$enDash = [char]0x2013
$body1 = "<w:body><w:p>" +
    '<w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change ' + $enDash + ' Ve</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r>' +
    "</w:p></w:body>"
$r1.InsertXML((New-BodyXml $body1))

# ---------------------------------------------------------------------------
# 2) Fourth paragraph ("Crispian's Day speech from Shakespear's Henry V
#    [Source – Wikipedia]"): keep "Crispian's" (with its proofErr wrapper)
#    untouched, restructure everything after it.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$afterCrispians = $d.Range($p4.Range.Start, $p4.Range.End)
$afterCrispians.Find.ClearFormatting()
$found = $afterCrispians.Find.Execute("Crispian" + [char]0x2019 + "s")
$restStart = $afterCrispians.End
$p4 = $d.Paragraphs.Item(4)
$restRange = $d.Range($restStart, $p4.Range.End - 1)

$rpr4 = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="202122"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
$body2 = "<w:body><w:p>" +
    '<w:r>' + $rpr4 + '<w:t xml:space="preserve"> Day speech from </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rpr4 + '<w:t>Shakespear' + [char]0x2019 + 's</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rpr4 + '<w:t xml:space="preserve"> Henry V [Source ' + $enDash + ' Wikipedia]</w:t></w:r>' +
    "</w:p></w:body>"
$restRange.InsertXML((New-BodyXml $body2))

# ---------------------------------------------------------------------------
# 3) Append two new paragraphs at the very end of the body (before sectPr):
#    one carrying the "larger" style plus shading/spacing overrides, and one
#    completely empty paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$body3 = "<w:body>" +
    '<w:p><w:pPr><w:pStyle w:val="larger"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0"/></w:pPr></w:p>' +
    '<w:p/>' +
    "</w:body>"
$endRange.InsertXML((New-BodyXml $body3))

# ---------------------------------------------------------------------------
# 4) Styles part: drop the now-unused "Hyperlink" and "apple-converted-space"
#    character styles (delete in reverse-index order to dodge a runtime bug
#    when re-looking-up a later style by name after an earlier one shifted).
# ---------------------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.Delete()
$appleStyle = $d.Styles.Item("apple-converted-space")
$appleStyle.Delete()
